$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- e100 (row 61): refresh body text, keep label/image name, resize image ---
$ws.Range("B61").Value = @'
<Bold>e100 Evening Debriefing</Bold> 
<LineBreak/><LineBreak/>
An evening debriefing is performed per <InlineUIContainer><Button Content='r4.9' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. 
Click image to continue to continue to Crew Rating Improvements per 
<InlineUIContainer><Button Content='r4.91' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<LineBreak/><LineBreak/>
             <InlineUIContainer><Image Name='Debrief' Height='200' Width='400'></Image></InlineUIContainer>
'@

# --- e101 (row 62): add rule references, fix typo, drop trailing paragraph ---
$ws.Range("A62").Value = "e101"
$ws.Range("B62").Value = @'
<Bold>e101 Evening Debriefing - Victory Point Total</Bold> 
<InlineUIContainer><Button Content='r4.92' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<InlineUIContainer><Button Content='r6.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<LineBreak/><LineBreak/>
The After Action Report 
<InlineUIContainer><Button Content='AAR' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
is updated to reflect victory points for both your tank and friendly forces. 
<LineBreak/><LineBreak/>
'@

# --- e102 (row 63): promotions now carries the full ladder + tracking fields ---
$ws.Range("A63").Value = "e102"
$ws.Range("B63").Value = @'
<Bold>e102 Evening Debriefing - Promotions</Bold> 
<InlineUIContainer><Button Content='r4.93' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<InlineUIContainer><Button Content='r25.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<LineBreak/><LineBreak/>
You are promoted when promotion points reach these values. However, cannot be promoted faster than once per month:
<LineBreak/><LineBreak/>
100 = Staff Sergeant<LineBreak/>
200 = 2nd Lieutenant<LineBreak/>
300 = 1st Lieutenant<LineBreak/>
400 = Captian
<LineBreak/><LineBreak/>
Promotion Points:  PROMOTION_POINTS<LineBreak/>
Promotion Date:     PROMOTION_DATE<LineBreak/>
<LineBreak/><LineBreak/>
'@

# --- e103 (row 64): decorations text expanded with roll target + modifiers header ---
$ws.Range("A64").Value = "e103"
$ws.Range("B64").Value = @'
<Bold>e103 Evening Debriefing - Decorations</Bold> 
<InlineUIContainer><Button Content='r26.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<LineBreak/><LineBreak/>
Roll for possible decorations on the 
<InlineUIContainer><Button Content='Decorations' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Table. 
Medals received are recorded on the After Action Report 
<InlineUIContainer><Button Content='AAR' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. Must roll at least 200 after modifications with 2D.
<LineBreak/><LineBreak/>
<Underline>Modifiers:</Underline><LineBreak/>
'@

# Row heights grew to fit the new copy
$ws.Rows(62).RowHeight = 120
$ws.Rows(63).RowHeight = 210
$ws.Rows(64).RowHeight = 135

# Scroll the view down one row, matching the saved window state (selection stays on B63)
try {
    $excel.ActiveWindow.ScrollRow = 62
} catch {
}
